$d = $word.ActiveDocument

# Paragraph 1: replace the placeholder token text and drop the trailing
# space-only run, update the left indent, and add a paragraph border.
$p1 = $d.Paragraphs(1)

$d.Content.Find.Execute(
    "**ID__AFFARS_pgi_5341_topic_2__ID** ", $true, $false, $false, $false,
    $false, $true, 1, $false, "**ID__AFFARS_AF_PGI_5341__ID**", 2) | Out-Null

$p1.LeftIndent = 11.25

$p1.Borders.DistanceFromTop = 5
$p1.Borders.DistanceFromLeft = 5
$p1.Borders.DistanceFromBottom = 5
$p1.Borders.DistanceFromRight = 5
